# design_input_background.xlsx edit:
#  - Select the cell that will become the active selection on "designinput"
#    (B10) BEFORE switching sheets, so the final ActiveSheet/activeTab ends
#    up on the newly-inserted "corr0" sheet (matches workbookView activeTab=3).
#  - Update PARAM5's numreal (B9) 10 -> 500.
#  - Point PARAM5 (row 9) and PARAM6 (row 10) at the new "corr0" correlation
#    sheet via the corr_sheet column (O9, O10).
#  - Insert a new worksheet "corr0" right before "corr1" containing the
#    2x2 correlation matrix between PARAM5 and PARAM6.

$wb = $excel.ActiveWorkbook

$designinput = $wb.Worksheets.Item("designinput")
$designinput.Range("B9").Value = 500
$designinput.Range("O9").Value = "corr0"
$designinput.Range("O10").Value = "corr0"
[void]$designinput.Range("B10").Select()

$corr1 = $wb.Worksheets.Item("corr1")
$corr0 = $wb.Worksheets.Add($corr1)
$corr0.Name = "corr0"

$corr0.Range("B1").Value = "PARAM5"
$corr0.Range("C1").Value = "PARAM6"
$corr0.Range("A2").Value = "PARAM5"
$corr0.Range("B2").Value = 1
$corr0.Range("A3").Value = "PARAM6"
$corr0.Range("B3").Value = 0.8
$corr0.Range("C3").Value = 1

[void]$corr0.Range("C8").Select()
